$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.189483029958524
$ws.Range("C2").Value = 0.04738906652589492
$ws.Range("E2").Value = 0.0736173268394662
$ws.Range("F2").Value = 3.232318159589937
$ws.Range("G2").Value = 0.002596794289770324
$ws.Range("I2").Value = 2.059841795425797
$ws.Range("J2").Value = 0.1688416664118542
$ws.Range("K2").Value = 0.9787752588754586
$ws.Range("L2").Value = 0.3249379500940179
$ws.Range("M2").Value = 0.3203958102452553
$ws.Range("N2").Value = 3.558174591686011

$ws.Range("B3").Value = 1.15796167866759
$ws.Range("C3").Value = 0.04216601487478044
$ws.Range("E3").Value = 0.07364045575576039
$ws.Range("F3").Value = 3.22648453291049
$ws.Range("G3").Value = 0.002600551967295604
$ws.Range("I3").Value = 2.063353970474147
$ws.Range("J3").Value = 0.1692261346921278
$ws.Range("K3").Value = 0.9439639254313192
$ws.Range("L3").Value = 0.3223180992656935
$ws.Range("M3").Value = 0.3144656596694198
$ws.Range("N3").Value = 3.577458116180772

$ws.Range("B4").Value = 1.139243508251411
$ws.Range("C4").Value = 0.0389763259432101
$ws.Range("E4").Value = 0.07367131596500087
$ws.Range("F4").Value = 3.224198904527071
$ws.Range("G4").Value = 0.002602982452954895
$ws.Range("I4").Value = 2.066241062877985
$ws.Range("J4").Value = 0.1694691947591052
$ws.Range("K4").Value = 0.9231443203630079
$ws.Range("L4").Value = 0.3208449698894427
$ws.Range("M4").Value = 0.3109842323870922
$ws.Range("N4").Value = 3.590137767256664

$ws.Range("B5").Value = 1.131775976078757
$ws.Range("C5").Value = 0.03768075663448656
$ws.Range("E5").Value = 0.07368809907461138
$ws.Range("F5").Value = 3.223593717063153
$ws.Range("G5").Value = 0.002604003986001116
$ws.Range("I5").Value = 2.067601368801874
$ws.Range("J5").Value = 0.1695700089675549
$ws.Range("K5").Value = 0.9147997706300828
$ws.Range("L5").Value = 0.3202788282346987
$ws.Range("M5").Value = 0.3096057833099444
$ws.Range("N5").Value = 3.59551593496235

$ws.Range("B6").Value = 1.130545688036619
$ws.Range("C6").Value = 0.03746588264284867
$ws.Range("E6").Value = 0.07369114054868398
$ws.Range("F6").Value = 3.223512937464875
$ws.Range("G6").Value = 0.002604175491187844
$ws.Range("I6").Value = 2.067838349468587
$ws.Range("J6").Value = 0.1695868558937317
$ws.Range("K6").Value = 0.913422600334286
$ws.Range("L6").Value = 0.3201868877748595
$ws.Range("M6").Value = 0.3093793275243115
$ws.Range("N6").Value = 3.596421726324493

$ws.Range("B7").Value = 1.139142149118811
$ws.Range("C7").Value = 0.03895883633661867
$ws.Range("E7").Value = 0.07367152524800513
$ws.Range("F7").Value = 3.224189421528166
$ws.Range("G7").Value = 0.00260299610379695
$ws.Range("I7").Value = 2.066258664175756
$ws.Range("J7").Value = 0.1694705472192819
$ws.Range("K7").Value = 0.9230312174097435
$ws.Range("L7").Value = 0.3208371962086503
$ws.Range("M7").Value = 0.3109654789923439
$ws.Range("N7").Value = 3.590209444340957

$ws.Range("B8").Value = 1.178482636573051
$ws.Range("C8").Value = 0.04558451373763717
$ws.Range("E8").Value = 0.07362185497687435
$ws.Range("F8").Value = 3.230037844960904
$ws.Range("G8").Value = 0.002598064413168236
$ws.Range("I8").Value = 2.060901199097522
$ws.Range("J8").Value = 0.1689727830956187
$ws.Range("K8").Value = 0.9666572455122093
$ws.Range("L8").Value = 0.3240065705064623
$ws.Range("M8").Value = 0.318318004231859
$ws.Range("N8").Value = 3.564649321609096

$ws.Range("B9").Value = 1.260667190518063
$ws.Range("C9").Value = 0.05871938903169394
$ws.Range("E9").Value = 0.07365586808071534
$ws.Range("F9").Value = 3.251783662762676
$ws.Range("G9").Value = 0.002589366962884475
$ws.Range("I9").Value = 2.056190347762723
$ws.Range("J9").Value = 0.1680518670339923
$ws.Range("K9").Value = 1.056609317450153
$ws.Range("L9").Value = 0.3312926063999271
$ws.Range("M9").Value = 0.3340002212908644
$ws.Range("N9").Value = 3.521186367328113

$ws.Range("B10").Value = 1.324115723942214
$ws.Range("C10").Value = 0.06846366451142671
$ws.Range("E10").Value = 0.07376000075233335
$ws.Range("F10").Value = 3.274020676924039
$ws.Range("G10").Value = 0.002583564294726887
$ws.Range("I10").Value = 2.056261801238385
$ws.Range("J10").Value = 0.1674084710454435
$ws.Range("K10").Value = 1.125389498153169
$ws.Range("L10").Value = 0.3372941079723972
$ws.Range("M10").Value = 0.3462897476069244
$ws.Range("N10").Value = 3.493313479448403

$ws.Range("B11").Value = 1.353645975845325
$ws.Range("C11").Value = 0.07291870012041102
$ws.Range("E11").Value = 0.0738243580552318
$ws.Range("F11").Value = 3.285495692021115
$ws.Range("G11").Value = 0.002581050733931255
$ws.Range("I11").Value = 2.057061489136359
$ws.Range("J11").Value = 0.1671228870898833
$ws.Range("K11").Value = 1.157266543000361
$ws.Range("L11").Value = 0.3401642862137209
$ws.Range("M11").Value = 0.3520468313041931
$ws.Range("N11").Value = 3.48151470586447

$ws.Range("B12").Value = 1.36492406447752
$ws.Range("C12").Value = 0.07460903096180971
$ws.Range("E12").Value = 0.07385115380931673
$ws.Range("F12").Value = 3.290036292144237
$ws.Range("G12").Value = 0.002580116947846428
$ws.Range("I12").Value = 2.057474610403361
$ws.Range("J12").Value = 0.1670157583185805
$ws.Range("K12").Value = 1.169422216053107
$ws.Range("L12").Value = 0.3412711982014685
$ws.Range("M12").Value = 0.3542507607622127
$ws.Range("N12").Value = 3.477173488966244

$ws.Range("B13").Value = 1.362490878478525
$ws.Range("C13").Value = 0.07424483941402116
$ws.Range("E13").Value = 0.07384527528951956
$ws.Range("F13").Value = 3.289049711573639
$ws.Range("G13").Value = 0.002580317254358134
$ws.Range("I13").Value = 2.057380732207577
$ws.Range("J13").Value = 0.1670387853268878
$ws.Range("K13").Value = 1.166800516705592
$ws.Range("L13").Value = 0.3410319153565382
$ws.Range("M13").Value = 0.3537750460815943
$ws.Range("N13").Value = 3.478102810471526

$ws.Range("B14").Value = 1.354571916369991
$ws.Range("C14").Value = 0.07305769785418192
$ws.Range("E14").Value = 0.07382651405725404
$ws.Range("F14").Value = 3.285865337755581
$ws.Range("G14").Value = 0.002580973549818529
$ws.Range("I14").Value = 2.05709326648379
$ws.Range("J14").Value = 0.167114053215589
$ws.Range("K14").Value = 1.15826490409907
$ws.Range("L14").Value = 0.3402549515330833
$ws.Range("M14").Value = 0.3522276725409839
$ws.Range("N14").Value = 3.481155011197359

$ws.Range("B15").Value = 1.349733766966779
$ws.Range("C15").Value = 0.07233097272222722
$ws.Range("E15").Value = 0.07381533755115477
$ws.Range("F15").Value = 3.283940239623504
$ws.Range("G15").Value = 0.002581377897037604
$ws.Range("I15").Value = 2.056931548681249
$ws.Range("J15").Value = 0.1671602890858477
$ws.Range("K15").Value = 1.153047599700955
$ws.Range("L15").Value = 0.339781645206088
$ws.Range("M15").Value = 0.3512829642266126
$ws.Range("N15").Value = 3.483041077523836

$ws.Range("B16").Value = 1.322199250893618
$ws.Range("C16").Value = 0.06817297475286921
$ws.Range("E16").Value = 0.07375613515813129
$ws.Range("F16").Value = 3.27329809260371
$ws.Range("G16").Value = 0.002583731091724999
$ws.Range("I16").Value = 2.056224975401946
$ws.Range("J16").Value = 0.1674272770437648
$ws.Range("K16").Value = 1.123318092770745
$ws.Range("L16").Value = 0.3371093446918536
$ws.Range("M16").Value = 0.3459168509884734
$ws.Range("N16").Value = 3.494102285981057

$ws.Range("B17").Value = 1.305478374168956
$ws.Range("C17").Value = 0.06562796185737341
$ws.Range("E17").Value = 0.07372415454412717
$ws.Range("F17").Value = 3.267117480560998
$ws.Range("G17").Value = 0.0025852069350679
$ws.Range("I17").Value = 2.055987992344569
$ws.Range("J17").Value = 0.1675928801984989
$ws.Range("K17").Value = 1.105230676557539
$ws.Range("L17").Value = 0.3355057751890627
$ws.Range("M17").Value = 0.3426674932231322
$ws.Range("N17").Value = 3.501113609032515

$ws.Range("B18").Value = 1.295923786618374
$ws.Range("C18").Value = 0.064166232428434
$ws.Range("E18").Value = 0.07370736006894774
$ws.Range("F18").Value = 3.26369051028496
$ws.Range("G18").Value = 0.002586067675486192
$ws.Range("I18").Value = 2.055923898092416
$ws.Range("J18").Value = 0.167688799356009
$ws.Range("K18").Value = 1.094882681033425
$ws.Range("L18").Value = 0.3345966341687756
$ws.Range("M18").Value = 0.3408142277256232
$ws.Range("N18").Value = 3.505229231640484

$ws.Range("B19").Value = 1.292699566347125
$ws.Range("C19").Value = 0.06367167278696684
$ws.Range("E19").Value = 0.07370194907698568
$ws.Range("F19").Value = 3.262552179548038
$ws.Range("G19").Value = 0.00258636114972476
$ws.Range("I19").Value = 2.055914600140966
$ws.Range("J19").Value = 0.1677213909327531
$ws.Range("K19").Value = 1.09138854735059
$ws.Range("L19").Value = 0.3342910831492816
$ws.Range("M19").Value = 0.3401894397791168
$ws.Range("N19").Value = 3.506636945904916

$ws.Range("B20").Value = 1.307251840881577
$ws.Range("C20").Value = 0.06589866533749955
$ws.Range("E20").Value = 0.07372739347699131
$ws.Range("F20").Value = 3.267762175945336
$ws.Range("G20").Value = 0.002585048600643842
$ws.Range("I20").Value = 2.056005746079549
$ws.Range("J20").Value = 0.1675751823002116
$ws.Range("K20").Value = 1.107150379197833
$ws.Range("L20").Value = 0.3356751135622744
$ws.Range("M20").Value = 0.3430117708840825
$ws.Range("N20").Value = 3.50035866187126

$ws.Range("B21").Value = 1.356895314841097
$ws.Range("C21").Value = 0.0734062996243523
$ws.Range("E21").Value = 0.07383195901113027
$ws.Range("F21").Value = 3.286795367809006
$ws.Range("G21").Value = 0.002580780290733581
$ws.Range("I21").Value = 2.057174708883366
$ws.Range("J21").Value = 0.1670919176890138
$ws.Range("K21").Value = 1.16076972723576
$ws.Range("L21").Value = 0.3404826216443695
$ws.Range("M21").Value = 0.3526815272416499
$ws.Range("N21").Value = 3.480255066041281

$ws.Range("B22").Value = 1.389897413866549
$ws.Range("C22").Value = 0.07833228498100198
$ws.Range("E22").Value = 0.07391442552824046
$ws.Range("F22").Value = 3.300372697489379
$ws.Range("G22").Value = 0.002578095840687467
$ws.Range("I22").Value = 2.058581578915934
$ws.Range("J22").Value = 0.1667819940794537
$ws.Range("K22").Value = 1.19630571499917
$ws.Range("L22").Value = 0.3437413539215868
$ws.Range("M22").Value = 0.3591402244024735
$ws.Range("N22").Value = 3.467854876968758

$ws.Range("B23").Value = 1.372232691484726
$ws.Range("C23").Value = 0.0757013947908689
$ws.Range("E23").Value = 0.07386912477884877
$ws.Range("F23").Value = 3.293022156636695
$ws.Range("G23").Value = 0.002579518991510209
$ws.Range("I23").Value = 2.057771888888055
$ws.Range("J23").Value = 0.1669468663924745
$ws.Range("K23").Value = 1.177294448612656
$ws.Range("L23").Value = 0.3419914587920658
$ws.Range("M23").Value = 0.355680414256021
$ws.Range("N23").Value = 3.474405477518687

$ws.Range("B24").Value = 1.306449874435856
$ws.Range("C24").Value = 0.06577627584165668
$ws.Range("E24").Value = 0.07372592419709889
$ws.Range("F24").Value = 3.267470315575935
$ws.Range("G24").Value = 0.002585120145456176
$ws.Range("I24").Value = 2.05599749485809
$ws.Range("J24").Value = 0.167583181301671
$ws.Range("K24").Value = 1.106282323696774
$ws.Range("L24").Value = 0.3355985158927268
$ws.Range("M24").Value = 0.3428560767413131
$ws.Range("N24").Value = 3.5006997095673

$ws.Range("B25").Value = 1.237895140053865
$ws.Range("C25").Value = 0.05515008240587349
$ws.Range("E25").Value = 0.07363269855534149
$ws.Range("F25").Value = 3.244801122278346
$ws.Range("G25").Value = 0.002591616265748482
$ws.Range("I25").Value = 2.05684439559667
$ws.Range("J25").Value = 0.1680518670339923
$ws.Range("K25").Value = 1.031802753232853
$ws.Range("L25").Value = 0.3292073497161567
$ws.Range("M25").Value = 0.3296227374771092
$ws.Range("N25").Value = 3.532231328249551
